$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.607716798782349
$ws.Range("B1").Value = 2.956317901611328
$ws.Range("C1").Value = 6.265285491943359
$ws.Range("D1").Value = 2.12846827507019
$ws.Range("E1").Value = 0.9476034641265869
